$wb = $excel.ActiveWorkbook
$mapping = $wb.Worksheets.Item("Mapping")

# Correct the Norwegian spelling/labels in the Mapping sheet. Ark1's cells
# pull these values in via formulas (e.g. =Mapping!D$3), so updating the
# source text here causes the dependent cached formula values to refresh too.
$mapping.Range("D3").Value = "1 - Harmløst (sluff)"
$mapping.Range("D4").Value = "2 - Små"
$mapping.Range("D7").Value = "5 - Svært store"
$mapping.Range("D11").Value = "Naturlig utløst"
$mapping.Range("D16").Value = "Få bratte heng"
$mapping.Range("D19").Value = "De fleste bratte heng, Også i mindre bratt terreng"

$excel.CalculateFull()

# Update the selected cell in each sheet, keeping "Ark1" as the active tab.
$mapping.Range("D10").Select()

$ark1 = $wb.Worksheets.Item("Ark1")
$ark1.Select()
$ark1.Range("B11").Select()
